$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round all numeric values in B2:E13 to the nearest integer,
# so the workbook stores integer data instead of full-precision floats.
$range = $ws.Range("B2:E13")
foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = [Math]::Floor([double]$val + 0.5)
    }
}
